$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 1.14
$ws.Range("N2").Value = 5.5

# Row 6
$ws.Range("Q6").Value = 2.3
$ws.Range("R6").Value = 1.6

# Row 9
$ws.Range("G9").Value = 2.45
$ws.Range("H9").Value = 3.4
$ws.Range("I9").Value = 2.8
$ws.Range("J9").Value = 3.2
$ws.Range("L9").Value = 3.5
$ws.Range("O9").Value = 1.33
$ws.Range("P9").Value = 3.25
$ws.Range("Q9").Value = 2.08
$ws.Range("R9").Value = 1.73
$ws.Range("W9").Value = 8
$ws.Range("X9").Value = 12
$ws.Range("Y9").Value = 10
$ws.Range("Z9").Value = 23
$ws.Range("AC9").Value = 9
$ws.Range("AE9").Value = 15
$ws.Range("AH9").Value = 13
$ws.Range("AJ9").Value = 29
$ws.Range("AK9").Value = 23
$ws.Range("AM9").Value = 301
$ws.Range("AN9").Value = 4.5
$ws.Range("AO9").Value = 15
$ws.Range("AQ9").Value = 51
$ws.Range("AU9").Value = 8
$ws.Range("AW9").Value = 4.75
$ws.Range("AY9").Value = 26
